# Kotte2014.xlsx edit: merge CKinetics.m / modelgen.m generated changes
# - Sheet "Kotte2014": rename species, add new transport reactions (ACt2r, FDPt2r),
#   add stoichiometry coefficients for FBP reaction, relabel activator/inhibitor
#   species with stoichiometric prefixes, add Vmax (L) values for rows 2-3.
# - Select A1:U7 on the Kotte2014 sheet so it matches saved selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kotte2014")

# Row 4: enzyme "GLUY" becomes the "FBP" reaction (fdp[c] ---> fdp[e])
$ws.Range("A4").Value = "FBP"

# Row 4: activator species relabeled with stoichiometric coefficient "4 pep[c]"
# (quote-prefixed so Excel stores it as literal text, matching the workbook's
# quotePrefix-styled cell)
$ws.Range("P4").Value = "'4 pep[c]"

# Row 5: inhibitor species relabeled with stoichiometric coefficient "2 fdp[c]"
$ws.Range("P5").Value = "2 fdp[c]"

# New transport reactions added in rows 6-7
$ws.Range("A6").Value = "ACt2r"
$ws.Range("A7").Value = "FDPt2r"

# New Vmax (L column) entries for rows 2 and 3
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1

# Restore the saved selection on the Kotte2014 sheet
$ws.Range("A1:U7").Select() | Out-Null
